$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spieltabelle")

# Row 5: update the Input (K5) amount
$ws.Range("K5").Value = 1.5

# Note: the text-cell writes below are intentionally ordered to match the
# shared-string table layout of the target workbook (the underlying engine
# appends/reuses shared-string slots in the order cells are written).
$ws.Range("R5").Value = "VPIP immer noch zu hoch, eigene Continuationbet war Fehler. Range des Gegners wurdenicht beachtet."

# Row 6: fill in the previously-empty session row
$ws.Range("B6").Value = "Cashgame"
$ws.Range("C6").Value = "sc.ch"
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 0.8
$ws.Range("F6").Value = 0.01
$ws.Range("G6").Value = 45934
$ws.Range("H6").Value = 0.71180555555555558
$ws.Range("I6").Value = 45934
$ws.Range("J6").Value = 0.76041666666666663
$ws.Range("K6").Value = 1.5
$ws.Range("L6").Value = 2.87
$ws.Range("M6").Value = 2.87
$ws.Range("N6").Value = 3
$ws.Range("O6").Value = 0.35
$ws.Range("P6").Value = "1 Schotch"

$ws.Range("Q5").Value = "2 führten nicht zum Ziel"

$ws.Range("Q6").Value = "Immer noch zu viel preflop mit wenig Value gecallt, viel Kartenglück (ein Vierling, ein Fullhouse)"
$ws.Range("R6").Value = "Gegener waren leichter zu bluffen, C-Bet ist aber immer noch ein Problem (Ich bin zu passiv)"

# Update the selected cell in the sheet view to mirror the saved-file state
$ws.Range("A8").Select()
